$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.989.17"
$ws.Range("E2").Value = "  +1.95%  "

$ws.Range("D3").Value = "1.669.57"
$ws.Range("E3").Value = "  +2.86%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.03"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("E6").Value = "  +2.15%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +2.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0617"
$ws.Range("E9").Value = "  +1.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.11"
$ws.Range("E10").Value = "  +4.86%  "

$ws.Range("E11").Value = "  +4.50%  "

$ws.Range("D12").Value = "1.906.56"
$ws.Range("E12").Value = "  +2.99%  "

$ws.Range("D13").Value = "1.660.80"
$ws.Range("E13").Value = "  +2.35%  "

$ws.Range("E14").Value = "  +0.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "65.72"
$ws.Range("E15").Value = "  +2.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.520"
$ws.Range("E16").Value = "  +1.61%  "

$ws.Range("D17").Value = "27.023.82"
$ws.Range("E17").Value = "  +2.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "234.83"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("E19").Value = "  +1.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.72"
$ws.Range("E20").Value = "  -0.80%  "

$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("E22").Value = "  +3.39%  "

$ws.Range("E23").Value = "  +1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.23"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.57"
$ws.Range("E25").Value = "  -0.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.16"
$ws.Range("E26").Value = "  +1.33%  "

$ws.Range("E27").Value = "  +0.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.88"
$ws.Range("E28").Value = "  +1.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("E31").Value = "  +1.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("E32").Value = "  +1.91%  "

$ws.Range("D33").Value = "1.448.75"
$ws.Range("E33").Value = "  -4.75%  "

$ws.Range("E34").Value = "  +5.38%  "

$ws.Range("E35").Value = "  +5.44%  "

$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("E37").Value = "  +0.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.892"
$ws.Range("E38").Value = "  +7.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0169"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("E40").Value = "  +3.55%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("E42").Value = "  +11.05%  "

$ws.Range("E43").Value = "  +3.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.67"
$ws.Range("E44").Value = "  +4.59%  "

$ws.Range("D45").Value = "1.814.13"
$ws.Range("E45").Value = "  +2.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.780"
$ws.Range("E46").Value = "  +2.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.48"
$ws.Range("E47").Value = "  +0.89%  "

$ws.Range("E48").Value = "  +1.43%  "

$ws.Range("E49").Value = "  +4.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0509"
$ws.Range("E50").Value = "  +1.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.60"
$ws.Range("E51").Value = "  +1.04%  "
